# AutoCommit_13 июня 2024 г. 8:47:58_SibNout2023
# Fills in additional grade cells (column I / J) for several rows and
# extends the colour-scale conditional formatting to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Cells that already carry a style (s="2") just need a value typed in
# ---------------------------------------------------------------------
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 5

$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5

$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 5

$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 5

# ---------------------------------------------------------------------
# 2) New cells that reuse the existing "green box" style (same style as
#    I8 / I24, fill + thick border on the right edge of the table)
# ---------------------------------------------------------------------
$styleSrc = $ws.Range("I8")

foreach ($addr in @("J8", "I18", "J18", "I22", "J22", "I23")) {
    $dst = $ws.Range($addr)
    $styleSrc.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Value = 5
}

# ---------------------------------------------------------------------
# 3) New cells that reuse the existing style used by J17 (fill-less,
#    thick border on both edges)
# ---------------------------------------------------------------------
$styleSrc2 = $ws.Range("J17")
$dst = $ws.Range("J28")
$styleSrc2.Copy()
$dst.PasteSpecial($xlPasteFormats)
$dst.Value = 5

# ---------------------------------------------------------------------
# 4) Brand new styles needed for I25/I26 (green fill + left-thick border,
#    same border as column J's outer edge) and J25/J26 (green fill, no
#    border). Build them from existing cells that already carry the
#    matching border, then paint the green fill on top.
# ---------------------------------------------------------------------
$greenFill = 5296274   # RGB(146,208,80) == theme green used elsewhere (FF92D050)

$borderSrc = $ws.Range("J2")   # already has the thick left border (borderId=2)
foreach ($addr in @("I25", "I26")) {
    $dst = $ws.Range($addr)
    $borderSrc.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Interior.Color = $greenFill
    $dst.Font.Bold = $false
    $dst.HorizontalAlignment = $xlCenter
    $dst.VerticalAlignment = $xlCenter
    $dst.WrapText = $true
    $dst.Value = 5
}

$noBorderSrc = $ws.Range("A1")   # plain, no border at all (borderId=0)
foreach ($addr in @("J25", "J26")) {
    $dst = $ws.Range($addr)
    $noBorderSrc.Copy()
    $dst.PasteSpecial($xlPasteFormats)
    $dst.Interior.Color = $greenFill
    $dst.HorizontalAlignment = $xlCenter
    $dst.VerticalAlignment = $xlCenter
    $dst.WrapText = $true
    $dst.Value = 5
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5) Conditional formatting: the colour-scale now needs to cover the
#    newly filled cells as well.
# ---------------------------------------------------------------------
foreach ($addr in @("J8", "J18", "J22", "J24:J26")) {
    $ws.Range($addr).FormatConditions.AddColorScale(3) | Out-Null
}

# ---------------------------------------------------------------------
# 6) Move the active selection the way the author left it.
# ---------------------------------------------------------------------
$ws.Range("J18").Select()

Write-Output "edit applied"
